# "menage dans la branche divers..." -- the deck currently has no slides
# at all (no <p:sldIdLst> in presentation.xml). Add the single title
# slide that the published deck ships with.

$p = $ppt.ActivePresentation

# ppLayoutTitle = 1 -> the slide master's first layout (slideLayout1.xml),
# which supplies the ctrTitle / subTitle placeholders used below. Adding
# the slide also creates the missing <p:sldIdLst> entry in presentation.xml.
$s = $p.Slides.Add(1, 1)

# Title placeholder (shape 1 on the "ctrTitle" layout) gets the headline;
# the subtitle placeholder (shape 2) is left empty, same as the target deck.
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Hello world!"
